$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)

# New data entered on Sheet3 (Gmail login: Email / password)
$ws3.Range("A1").Value = "Email"
$ws3.Range("B1").Value = "password"
$ws3.Range("A2").Value = "charmproject17@gmail.com"
$ws3.Range("B2").Value = "Charm@1234"

# Excel auto-linkifies the typed e-mail / text -> Hyperlink style gets applied
$ws3.Hyperlinks.Add($ws3.Range("A2"), "mailto:charmproject17@gmail.com")
$ws3.Hyperlinks.Add($ws3.Range("B2"), "Charm@1234")

# Columns best-fit to their new content
$ws3.Columns("A:B").AutoFit() | Out-Null

# Sheet3 becomes the active sheet / tab
$ws3.Activate() | Out-Null
$ws3.Range("B6").Select() | Out-Null

Write-Host "done"
